# chore: update Sheets via scheduled runner
# Applies the scraped profit-recalculation edits across the ALC/ARM/BSM/CRP/
# CUL/GSM/LTW/WVR sheets (source file was named Sheets/Sophia_Profits.xlsx).

$wb = $excel.ActiveWorkbook

function Set-Cells($SheetName, $Cells) {
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($addr in $Cells.Keys) {
        $ws.Range($addr).Value = $Cells[$addr]
    }
}

# ---- ALC ----
Set-Cells "ALC" @{
    "H9"  = 200098
    "I9"  = 333413.34
    "K9"  = 333413.34
    "M9"  = -333244.34

    "H52" = 4500
    "I52" = 4500
    "K52" = 13500
    "M52" = -13340

    "H76" = 1999
    "I76" = 1999
    "K76" = 1999
    "M76" = -1684

    "H79" = 1999
    "I79" = 1999
    "K79" = 1999
    "M79" = -907

    "H100" = 6999.6665
    "I100" = 9999.5
    "J100" = 1000
    "K100" = 9999.5
    "L100" = 1000
    "M100" = -9458.5
    "N100" = -2082

    "H112" = 2482.8572
    "J112" = 2485.05
    "L112" = 7455.150000000001
    "N112" = -9671.150000000001

    "H141" = 5545.6665
    "I141" = 7998.5
    "J141" = 640
    "K141" = 23995.5
    "L141" = 1920
    "M141" = -18815.5
    "N141" = -12280
}

# ---- ARM ----
Set-Cells "ARM" @{
    "H32" = 5776.976
    "I32" = 5734.9756
    "K32" = 5734.9756
    "M32" = -5447.9756

    "H46" = 19990
    "J46" = 19990
    "L46" = 19990
    "N46" = -20628

    "H61" = 2777.6667
    "I61" = 2777.6667
    "K61" = 2777.6667
    "M61" = -2565.6667

    "H74" = 6722.613
    "I74" = 6455.6665
    "K74" = 6455.6665
    "M74" = -5581.6665

    "H77" = 6722.613
    "I77" = 6455.6665
    "K77" = 32278.3325
    "M77" = -27910.3325

    "H136" = 2777.6667
    "I136" = 2777.6667
    "K136" = 8333.000100000001
    "M136" = -5783.000100000001
}

# ---- BSM ----
Set-Cells "BSM" @{
    "H107" = 749.5
    "I107" = 749.5
    "K107" = 749.5
    "M107" = 1170.5
}

# ---- CRP ----
Set-Cells "CRP" @{
    "H31" = 1110.5
    "I31" = 1110.5
    "K31" = 1110.5
    "M31" = -815.5

    "H34" = 1110.5
    "I34" = 1110.5
    "K34" = 1110.5
    "M34" = -908.5

    "H58" = 3512.625
    "I58" = 3500.3572
    "K58" = 3500.3572
    "M58" = -3297.3572

    "H134" = 6930.143
    "I134" = 7685.3335
    "K134" = 23056.0005
    "M134" = -20521.0005

    "H136" = 3512.625
    "I136" = 3500.3572
    "K136" = 10501.0716
    "M136" = -7951.071599999999
}

# ---- CUL ----
Set-Cells "CUL" @{
    "H2" = 139.8
    "I2" = 98
    "J2" = 144.44444
    "K2" = 588
    "L2" = 866.6666399999999
    "M2" = -475
    "N2" = -1092.66664

    "H74" = 7253.3335
    "J74" = 7678.5713
    "L74" = 23035.7139
    "N74" = -25157.7139

    "H77" = 7253.3335
    "J77" = 7678.5713
    "L77" = 69107.14169999999
    "N77" = -79715.14169999999

    "H86" = 997.25
    "I86" = 932
    "J86" = 1019
    "K86" = 2796
    "L86" = 3057
    "M86" = -1610
    "N86" = -5429

    "H89" = 997.25
    "I89" = 932
    "J89" = 1019
    "K89" = 8388
    "L89" = 9171
    "M89" = -2460
    "N89" = -21027

    "H106" = 0
    "J106" = 0
    "L106" = 0

    "H113" = 1349
    "I113" = 331
    "J113" = 2112.5
    "K113" = 993
    "L113" = 6337.5
    "M113" = 1177
    "N113" = -10677.5

    "H137" = 994
    "I137" = 0
    "J137" = 994
    "K137" = 0
    "L137" = 2982
    "N137" = -13182
}
# row 106's N cell (-16892) and row 137's M cell (2115) are dropped entirely
$wsCUL = $wb.Worksheets.Item("CUL")
$wsCUL.Range("N106").ClearContents()
$wsCUL.Range("M137").ClearContents()

# ---- GSM ----
Set-Cells "GSM" @{
    "H107" = 2899.5
    "I107" = 2899.5
    "J107" = 0
    "K107" = 2899.5
    "L107" = 0
    "M107" = -979.5

    "H132" = 3728.7144
    "I132" = 2800.75
    "J132" = 4966
    "K132" = 8402.25
    "L132" = 14898
    "M132" = -5872.25
    "N132" = -19958
}
# row 107's N cell (-6640) is dropped entirely
$wsGSM = $wb.Worksheets.Item("GSM")
$wsGSM.Range("N107").ClearContents()

# ---- LTW ----
Set-Cells "LTW" @{
    "H7" = 598.2
    "I7" = 496.5
    "K7" = 496.5
    "M7" = -384.5

    "H126" = 598.2
    "I126" = 496.5
    "K126" = 1489.5
    "M126" = 980.5

    "H132" = 3504.3333
    "I132" = 1593.5883
    "K132" = 4780.7649
    "M132" = -2250.7649
}

# ---- WVR ----
Set-Cells "WVR" @{
    "H126" = 3125
    "I126" = 3166.6667
    "K126" = 9500.000100000001
    "M126" = -7030.000100000001

    "H132" = 1539.2354
    "I132" = 536.2308
    "K132" = 1608.6924
    "M132" = 921.3075999999999

    "H136" = 2713.5833
    "I136" = 2896.476
    "J136" = 1433.3334
    "K136" = 8689.428
    "L136" = 4300.0002
    "M136" = -6139.428
    "N136" = -9400.0002
}
